# "updated GME & AMC"
#
# The "Games" sub-list (E16:E18) under the E15 "Games" header is being
# reorganized:
#   - E16 "Blackjack"  -> "Card Games (Non-Poker)"
#   - E17 "Chess"      -> unchanged text, but becomes a hyperlinked entry
#   - E18 "Poker"      -> "Interactive Entertainment" (hyperlinked)
#   - E19 (new row)    -> "Poker" (the category that used to live in E18)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new/changed text values. E18 is written before E16 so the
# shared-string table picks up "Interactive Entertainment" then
# "Card Games (Non-Poker)" in that order (matches the source order).
$ws.Range("E18").Value = "Interactive Entertainment"
$ws.Range("E16").Value = "Card Games (Non-Poker)"

# New row: Poker moves down to E19.
$ws.Range("E19").Value = "Poker"

# Add the new hyperlinks (E18 first, then E17) so they match the
# relationship-id ordering of the target file.
$ws.Hyperlinks.Add($ws.Range("E18"), "..\..\OneDrive\Documents\Interactive Entertainment.xlsx")
$ws.Hyperlinks.Add($ws.Range("E17"), "..\..\OneDrive\Documents\Chess.xlsx")

# Give both newly-linked cells the standard "Hyperlink" cell style
# (applied after Hyperlinks.Add so the cells reuse the existing
# Hyperlink style entry instead of a fresh duplicate).
$ws.Range("E18").Style = "Hyperlink"
$ws.Range("E17").Style = "Hyperlink"

# Move the active selection to E15, matching the saved cursor position.
$ws.Range("E15").Select()
